# Update the "弹幕数" (danmu/comment count, column F) values for a handful of
# rows across the 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types)
# sheets, reflecting the regenerated data snapshot referenced by the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> list of (row, newValue) updates to column F.
$updates = @{
    "展览" = @(
        @{ Row = 9;  Value = 4772 },
        @{ Row = 17; Value = 2992 },
        @{ Row = 18; Value = 1825 },
        @{ Row = 25; Value = 310 },
        @{ Row = 27; Value = 2845 },
        @{ Row = 29; Value = 2491 },
        @{ Row = 37; Value = 934 },
        @{ Row = 42; Value = 204 },
        @{ Row = 45; Value = 3496 }
    )
    "演出" = @(
        @{ Row = 5; Value = 1 }
    )
    "全部类型" = @(
        @{ Row = 10; Value = 4772 },
        @{ Row = 15; Value = 2992 },
        @{ Row = 17; Value = 1825 },
        @{ Row = 27; Value = 310 },
        @{ Row = 28; Value = 2845 },
        @{ Row = 32; Value = 2491 },
        @{ Row = 39; Value = 934 },
        @{ Row = 49; Value = 3496 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Cells.Item($entry.Row, 6).Value = $entry.Value
    }
}
